# feat: adiciona suporte a expand_dates_to para transformar colunas de
# data em Ano, Mes e Valor
#
# Expands the date column into two extra columns (F, G) on the "Base"
# sheet: new headers in row 1, numeric values below, matching the
# formatting already used by the existing header/body cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Base")

# Copy the existing header cell formatting onto the two new header cells
# so they share the same style as the rest of row 1.
$ws.Range("E1").Copy() | Out-Null
$ws.Range("F1:G1").PasteSpecial(-4122) | Out-Null # xlPasteFormats

# Copy the existing body cell formatting onto the new data cells.
$ws.Range("D2").Copy() | Out-Null
$ws.Range("F2:G4").PasteSpecial(-4122) | Out-Null # xlPasteFormats

$ws.Application.CutCopyMode = 0

# New header labels
$ws.Range("F1").Value = "30/11/20"
$ws.Range("G1").Value = "30/04/20"

# New numeric values for rows 2-4
$ws.Range("F2").Value = 0.85
$ws.Range("G2").Value = 0.85

$ws.Range("F3").Value = 0.85
$ws.Range("G3").Value = 0.08

$ws.Range("F4").Value = 0.85
$ws.Range("G4").Value = 0.6
